$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: System 4 Truth State (column E) corrected from 1 to 2
$ws.Range("E8").Value = 2

# Row 10: System 2 Truth State (column C) corrected from 1 to 2
$ws.Range("C10").Value = 2

# Rows 16-24: System 1 (B) and System 3 (D) Truth States corrected from 2 to 1
for ($r = 16; $r -le 24; $r++) {
    $ws.Cells.Item($r, 2).Value = 1   # column B
    $ws.Cells.Item($r, 4).Value = 1   # column D
}

# Rows 19-25: System 4 Truth State (column F) corrected from 2 to 0
for ($r = 19; $r -le 25; $r++) {
    $ws.Cells.Item($r, 6).Value = 0   # column F
}

# Rows 45-102: System 2 Truth State (column C) corrected from 1 to 0
for ($r = 45; $r -le 102; $r++) {
    $ws.Cells.Item($r, 3).Value = 0   # column C
}
